# Update specific values in column B (per diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3   = 5.698800000000002
    14  = 5.917700000000001
    21  = 9.384200000000002
    23  = 8.893599999999996
    25  = 5.421700000000002
    26  = 4.624800000000003
    29  = 5.102300000000003
    53  = 5.629299999999997
    57  = 4.799999999999994
    59  = 5.1753
    69  = 5.407999999999992
    79  = 9.521600000000005
    83  = 5.337399999999997
    91  = 5.732800000000001
    93  = 5.454300000000001
    103 = 5.865699999999996
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
